$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24: SQL Saturday St. Louis 2025 (#1117) ---
$ws.Hyperlinks.Add($ws.Range("B24"), "https://sqlsaturday.com/2025-10-25-sqlsaturday1117/", "", "", "https://sqlsaturday.com/2025-10-25-sqlsaturday1117/")
$ws.Range("A24").Value = "10/25/2025"
$ws.Range("B24").Value = "SQL Saturday St. Louis 2025 (#1117)"
$ws.Range("B24").Style = "Hyperlink"
$ws.Range("C24").Value = 235
$ws.Range("D24").Value = 121
$ws.Range("E24").Formula = "=IF(C24=0,0,+(C24-D24)/C24)"

# --- Row 25: SQL Saturday Toronto 2025 (#1131) ---
$ws.Hyperlinks.Add($ws.Range("B25"), "https://sqlsaturday.com/2025-10-25-sqlsaturday1131/", "", "", "https://sqlsaturday.com/2025-10-25-sqlsaturday1131/")
$ws.Range("A25").Value = "10/25/2025"
$ws.Range("B25").Value = "SQL Saturday Toronto 2025 (#1131)"
$ws.Range("B25").Style = "Hyperlink"
$ws.Range("C25").Value = 367
$ws.Range("D25").Value = 289
$ws.Range("E25").Formula = "=IF(C25=0,0,+(C25-D25)/C25)"

# --- Row 23: SQL Saturday Pittsburgh 2025 (#1123) ---
$ws.Hyperlinks.Add($ws.Range("B23"), "https://sqlsaturday.com/2025-10-18-sqlsaturday1123/", "", "", "https://sqlsaturday.com/2025-10-18-sqlsaturday1123/")
$ws.Range("A23").Value = "10/18/2025"
$ws.Range("B23").Value = "SQL Saturday Pittsburgh 2025 (#1123) "
$ws.Range("B23").Style = "Hyperlink"
$ws.Range("E23").Formula = "=IF(C23=0,0,+(C23-D23)/C23)"

# --- Update selection to match final state ---
$ws.Range("D16").Select()
